$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Weekly refresh: a new reporting week's data was prepended to the
# "Vega Monumental Concepción - Acelga" series, which pushes every
# existing record down by two rows (one new "Primera"/"Segunda" pair).
#
# Step 1: open up two blank rows right before the block that needs to
# shift (rows 234-293 slide down to 236-295; rows 294/295 appear for
# the first time, carrying what used to be the last two records).
$ws.Range("A234:A235").EntireRow.Insert()

# Step 2: the two freshly-inserted blank rows should hold exactly the
# data that row 232/233 had *before* this week's update - copy it down.
$ws.Range("A232:R232").Copy($ws.Range("A234:R234"))
$ws.Range("A233:R233").Copy($ws.Range("A235:R235"))

# Step 3: now overwrite row 232/233 (the "Primera"/"Segunda" pair for
# this market) with this week's fresh figures.
$ws.Range("D232").Value = 44798
$ws.Range("K232").Value = 700
$ws.Range("L232").Value = 800
$ws.Range("M232").Value = 750
$ws.Range("P232").Value = 750

$ws.Range("D233").Value = 44798
$ws.Range("K233").Value = 600
$ws.Range("L233").Value = 600
$ws.Range("M233").Value = 600
$ws.Range("P233").Value = 600
